$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated K (formerly Strike#) column values per new simulation run
$kValues = @{
    2 = 3
    3 = 1
    4 = 1
    5 = 0
    6 = 1
    7 = 0
    8 = 1
    9 = 1
    10 = 4
    11 = 1
    12 = 2
    13 = 1
    14 = 2
    15 = 1
    16 = 0
    17 = 3
    18 = 0
    19 = 1
    20 = 0
    21 = 1
    22 = 0
    23 = 2
    24 = 0
    25 = 2
    26 = 2
    27 = 2
    28 = 0
    29 = 0
    30 = 1
    31 = 3
    32 = 1
    33 = 2
    34 = 2
    35 = 3
    36 = 0
    37 = 2
    38 = 0
    39 = 1
    40 = 0
    41 = 2
    42 = 0
    43 = 1
    44 = 0
    45 = 1
    46 = 2
    47 = 1
    48 = 1
    49 = 1
    50 = 2
    51 = 0
    52 = 0
    53 = 1
    54 = 2
    55 = 1
    56 = 1
    57 = 2
    58 = 1
    59 = 0
    60 = 2
    61 = 3
    62 = 2
    63 = 0
    64 = 1
    65 = 1
    66 = 0
    67 = 0
    68 = 1
    69 = 1
    70 = 0
    71 = 1
    72 = 4
    73 = 1
    74 = 0
    75 = 4
    76 = 1
    77 = 0
    78 = 3
    79 = 1
    80 = 0
    81 = 2
    82 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $kValues[$row]
}
